# Add a new "UK" market sheet, cloned from the existing "Poland" sheet,
# mirroring how the other per-country market sheets are laid out.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Poland")

# Copy "Poland" and place the copy immediately after it (i.e. as the new
# last sheet in the workbook) - this is exactly what Excel's own
# "Move or Copy... > Create a copy" does.
$source.Copy([System.Reflection.Missing]::Value, $source)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# Fill in the market-specific values (User Story / Jira ref first, then the
# market name) so new shared-string entries are appended in the same order
# as the authored workbook.
$newSheet.Range("B4").Value = "NGC-2741/T3343"
$newSheet.Range("B2").Value = "UK Market"

# Make the new UK sheet the active tab/selection, matching the saved state
# (this also naturally moves the "tabSelected" flag off of whatever sheet
# had it before - e.g. Norway).
$newSheet.Activate()
$newSheet.Range("B4").Select()
